{"js": "// Fix typo in report:\n//  1. \"Set up 5 matrices of size height x weight:\" -> \"...height x width:\"\n//     (Word also re-anchors its internal \"last edit\" (_GoBack) bookmark to\n//     sit right after the corrected word, so we mirror that.)\n//  2. \"A nice descr\" + \"i\" + \"ptive figure\" (three runs, same rendered text)\n//     collapse back down into a single \"A nice descriptive figure\" run once\n//     Word re-paints the hyperlink text.\n\nconst body = context.document.body;\n\n// --- 1. weight -> width -------------------------------------------------\nconst typo = body.search(\"weight\", { matchCase: true, matchWholeWord: true });\ntypo.load(\"text\");\nawait context.sync();\n\nif (typo.items.length > 0) {\n  const fixedRange = typo.items[0].insertText(\"width\", \"Replace\");\n  await context.sync();\n\n  // Drop the old \"_GoBack\" bookmark (it lived on an empty paragraph further\n  // down the document) before re-creating it at the new edit location --\n  // bookmark names must stay unique.\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n\n  fixedRange.getRange(\"End\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// --- 2. Re-merge the hyperlink's run-split text -------------------------\n// The three runs render identical text to one run, so search for the full\n// phrase (search spans run boundaries), swap to a placeholder and back to\n// force Word to collapse them into a single run, matching the final state.\nconst hyperlinkText = \"A nice descriptive figure\";\nconst split = body.search(hyperlinkText, { matchCase: true });\nsplit.load(\"text\");\nawait context.sync();\n\nif (split.items.length > 0) {\n  const placeholder = split.items[0].insertText(\"\\u0001TEMP_MERGE\\u0001\", \"Replace\");\n  await context.sync();\n\n  const placeholderHit = body.search(\"\\u0001TEMP_MERGE\\u0001\", { matchCase: true });\n  placeholderHit.load(\"text\");\n  await context.sync();\n\n  if (placeholderHit.items.length > 0) {\n    placeholderHit.items[0].insertText(hyperlinkText, \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Fix typo in report:\n#  1. \"Set up 5 matrices of size height x weight:\" -> \"...height x width:\"\n#     (Word also re-anchors its internal \"last edit\" (_GoBack) bookmark to\n#     sit right after the corrected word, so we mirror that.)\n#  2. \"A nice descr\" + \"i\" + \"ptive figure\" (three runs, same rendered text)\n#     collapse back down into a single \"A nice descriptive figure\" run once\n#     Word re-paints the hyperlink text.\n\n$d = $word.ActiveDocument\n\n# --- 1. weight -> width --------------------------------------------------\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"weight\", $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, \"width\", $wdReplaceAll)\n\nif ($found) {\n    # Drop the old \"_GoBack\" bookmark (it lived on an empty paragraph further\n    # down the document) before re-creating it at the new edit location --\n    # bookmark names must stay unique.\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n\n    $endRng = $rng.Duplicate\n    $endRng.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $endRng)\n}\n\n# --- 2. Re-merge the hyperlink's run-split text ---------------------------\n# The three runs render identical text to one run, so find the full phrase\n# (Find spans run boundaries), swap to a placeholder and back to force Word\n# to collapse them into a single run, matching the final state.\n$hyperlinkText = \"A nice descriptive figure\"\n$placeholder = \"TEMP_MERGE_PLACEHOLDER\"\n\n$rng2 = $d.Content\n$split = $rng2.Find.Execute($hyperlinkText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $placeholder, $wdReplaceAll)\n\nif ($split) {\n    $rng3 = $d.Content\n    $rng3.Find.Execute($placeholder, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $hyperlinkText, $wdReplaceAll) | Out-Null\n}\n"}
